# Apply the "Added more test data from excel file" change:
#  1. Append two new data rows (SQL, HTML) to the ShareSkill sheet.
#  2. Widen a couple of columns on the ShareSkill sheet to fit the new data.
#  3. Move the selection/active cell on the ShareSkill sheet.
#  4. Add a brand new "EditSkill" worksheet (after "DeleteSkill") with a
#     small title/credit-amount table, matching the DeleteSkill sheet layout.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Add the new "EditSkill" worksheet after "DeleteSkill" (its new
#    shared strings are registered first, matching the source order).
# ---------------------------------------------------------------------
$deleteSkill = $wb.Worksheets.Item("DeleteSkill")
$editSkill = $wb.Worksheets.Add($null, $deleteSkill)
$editSkill.Name = "EditSkill"

$editSkill.Range("A1").Value = "Title"
$editSkill.Range("B1").Value = "EditedCreditAmount"
$editSkill.Range("A2").Value = "Artificial Intelligence"
$editSkill.Range("B2").Value = 6

$editSkill.Cells.Item(1, 1).ColumnWidth = 19.333333333333332
$editSkill.Cells.Item(1, 2).ColumnWidth = 17.0

[void]$editSkill.Range("C10").Select()

# ---------------------------------------------------------------------
# 2. ShareSkill sheet updates: two new data rows (SQL, HTML)
# ---------------------------------------------------------------------
$share = $wb.Worksheets.Item("ShareSkill")

# New row 3: SQL
$share.Range("A3").Value = "SQL"
$share.Range("B3").Value = "Programming Language"
$share.Range("C3").Value = "Programming & Tech"
$share.Range("D3").Value = "Databases"
$share.Range("E3").Value = "AddingData"
$share.Range("F3").Value = "Hourly basis service"
$share.Range("G3").Value = "Online"

$share.Range("H2").Copy($share.Range("H3"))
$share.Range("H3").Value = 44638
$share.Range("I2").Copy($share.Range("I3"))
$share.Range("I3").Value = 44638

$share.Range("J3").Value = "Fri"
$share.Range("K3").Value = "10:00AM"
$share.Range("L3").Value = "11:00AM"
$share.Range("M3").Value = "Skill-Exchange"
$share.Range("N3").Value = "SQL Queries"
$share.Range("O3").Value = "Credit"
$share.Range("P3").Value = 4
$share.Range("Q3").Value = "Active"

# New row 4: HTML
$share.Range("A4").Value = "HTML"
$share.Range("B4").Value = "Markup Language"
$share.Range("C4").Value = "Programming & Tech"
$share.Range("D4").Value = "Web & Mobile App"
$share.Range("E4").Value = "Webpage Design"
$share.Range("F4").Value = "Hourly basis service"
$share.Range("G4").Value = "Online"

$share.Range("H2").Copy($share.Range("H4"))
$share.Range("H4").Value = 44641
$share.Range("I2").Copy($share.Range("I4"))
$share.Range("I4").Value = 44641

$share.Range("J4").Value = "Mon"
$share.Range("K4").Value = "8:00AM"
$share.Range("L4").Value = "9:00AM"
$share.Range("M4").Value = "Skill-Exchange"
$share.Range("N4").Value = "Web Design"
$share.Range("O4").Value = "Credit"
$share.Range("P4").Value = 6
$share.Range("Q4").Value = "Active"

# Widen column B and column D to fit the new content.
$share.Cells.Item(1, 2).ColumnWidth = 24.833333333333332
$share.Cells.Item(1, 4).ColumnWidth = 22.5

# Update view: drop the frozen "topLeftCell" scroll position and move the
# active selection to C14.
$share.Activate()
[void]$share.Range("C14").Select()
